$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    "Composition",
    "Bi2Ti2O7",
    "Bi1.5Al0.5Ti2O7",
    "Bi1.5Sc0.5Ti2O7",
    "Bi1.5Ga0.5Ti2O7",
    "Bi1.5Y0.5Ti2O7",
    "Bi1.5In0.5Ti2O7",
    "Bi1.5La0.5Ti2O7",
    "Bi4Ti3O12",
    "Bi3.5Al0.5Ti3O12",
    "Bi3.5Ga0.5Ti3O12",
    "Bi3.5In0.5Ti3O12",
    "Bi3.5Sc0.5Ti3O12",
    "Bi3.5Y0.5Ti3O12",
    "Bi3.5La0.5Ti3O12",
    "Bi4Ti2.5Al0.5O12",
    "Bi4Ti2.5Ga0.5O12",
    "Bi4Ti2.5Sc0.5O12",
    "Bi12TiO20",
    "BaTaO2N",
    "BaTa0.5Al0.5O2N",
    "BaTa0.5Mg0.5O2N",
    "BaTa0.5Al0.25Mg0.25O2N",
    "BaTa0.5Al0.125Mg0.375O2N",
    "BaTa0.5Al0.375Mg0.125O2N",
    "Al2O3",
    "TiO2",
    "Ta2O5",
    "Fe2O3",
    "NaLaTi2O6",
    "NaCl",
    "Bi1.5Cr0.5Ti2O7",
    "Bi1.5Fe0.5Ti2O7",
    "Bi4Ti2O11",
    "Ta1.94Eu0.06O5"
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $data[$i]
}

$ws.Range("A1:A35").Select() | Out-Null
